$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-10-16 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-17 Friday", 2) | Out-Null
$d.Content.Find.Execute("51+22=", $true, $false, $false, $false, $false, $true, 1, $false, "5+35=", 2) | Out-Null
$d.Content.Find.Execute("89-0=", $true, $false, $false, $false, $false, $true, 1, $false, "72-66=", 2) | Out-Null
$d.Content.Find.Execute("18+18=", $true, $false, $false, $false, $false, $true, 1, $false, "35+55=", 2) | Out-Null
$d.Content.Find.Execute("58-51=", $true, $false, $false, $false, $false, $true, 1, $false, "67-18=", 2) | Out-Null
$d.Content.Find.Execute("79-26=", $true, $false, $false, $false, $false, $true, 1, $false, "91-49=", 2) | Out-Null
$d.Content.Find.Execute("2+76=", $true, $false, $false, $false, $false, $true, 1, $false, "21+52=", 2) | Out-Null
$d.Content.Find.Execute("35-32=", $true, $false, $false, $false, $false, $true, 1, $false, "44-21=", 2) | Out-Null
$d.Content.Find.Execute("82+17=", $true, $false, $false, $false, $false, $true, 1, $false, "56-45=", 2) | Out-Null
$d.Content.Find.Execute("55-49=", $true, $false, $false, $false, $false, $true, 1, $false, "85-38=", 2) | Out-Null
$d.Content.Find.Execute("89-69=", $true, $false, $false, $false, $false, $true, 1, $false, "64-49=", 2) | Out-Null
$d.Content.Find.Execute("86-50=", $true, $false, $false, $false, $false, $true, 1, $false, "58-40=", 2) | Out-Null
$d.Content.Find.Execute("69-53=", $true, $false, $false, $false, $false, $true, 1, $false, "6+61=", 2) | Out-Null
$d.Content.Find.Execute("35+45=", $true, $false, $false, $false, $false, $true, 1, $false, "89-11=", 2) | Out-Null
$d.Content.Find.Execute("38+33=", $true, $false, $false, $false, $false, $true, 1, $false, "93-84=", 2) | Out-Null
$d.Content.Find.Execute("88-56=", $true, $false, $false, $false, $false, $true, 1, $false, "26+31=", 2) | Out-Null
$d.Content.Find.Execute("36-13=", $true, $false, $false, $false, $false, $true, 1, $false, "94-6=", 2) | Out-Null
$d.Content.Find.Execute("71-61=", $true, $false, $false, $false, $false, $true, 1, $false, "26-0=", 2) | Out-Null
$d.Content.Find.Execute("52-35=", $true, $false, $false, $false, $false, $true, 1, $false, "63-50=", 2) | Out-Null
$d.Content.Find.Execute("59-5=", $true, $false, $false, $false, $false, $true, 1, $false, "82-67=", 2) | Out-Null
$d.Content.Find.Execute("74-67=", $true, $false, $false, $false, $false, $true, 1, $false, "15+77=", 2) | Out-Null
$d.Content.Find.Execute("67-9=", $true, $false, $false, $false, $false, $true, 1, $false, "47-32=", 2) | Out-Null
$d.Content.Find.Execute("77-72=", $true, $false, $false, $false, $false, $true, 1, $false, "63+9=", 2) | Out-Null
$d.Content.Find.Execute("24+60=", $true, $false, $false, $false, $false, $true, 1, $false, "7+17=", 2) | Out-Null
$d.Content.Find.Execute("50-34=", $true, $false, $false, $false, $false, $true, 1, $false, "60-25=", 2) | Out-Null
$d.Content.Find.Execute("77-71=", $true, $false, $false, $false, $false, $true, 1, $false, "91+6=", 2) | Out-Null
$d.Content.Find.Execute("81-29=", $true, $false, $false, $false, $false, $true, 1, $false, "24+20=", 2) | Out-Null
$d.Content.Find.Execute("35+0=", $true, $false, $false, $false, $false, $true, 1, $false, "32-2=", 2) | Out-Null
$d.Content.Find.Execute("2+54=", $true, $false, $false, $false, $false, $true, 1, $false, "35-11=", 2) | Out-Null
$d.Content.Find.Execute("66+19=", $true, $false, $false, $false, $false, $true, 1, $false, "88-14=", 2) | Out-Null
$d.Content.Find.Execute("1+31=", $true, $false, $false, $false, $false, $true, 1, $false, "38-18=", 2) | Out-Null
$d.Content.Find.Execute("21-21=", $true, $false, $false, $false, $false, $true, 1, $false, "6+21=", 2) | Out-Null
$d.Content.Find.Execute("20+56=", $true, $false, $false, $false, $false, $true, 1, $false, "40+3=", 2) | Out-Null
$d.Content.Find.Execute("77-54=", $true, $false, $false, $false, $false, $true, 1, $false, "16+78=", 2) | Out-Null
$d.Content.Find.Execute("16-12=", $true, $false, $false, $false, $false, $true, 1, $false, "58-41=", 2) | Out-Null
$d.Content.Find.Execute("60-9=", $true, $false, $false, $false, $false, $true, 1, $false, "16+44=", 2) | Out-Null
$d.Content.Find.Execute("47+44=", $true, $false, $false, $false, $false, $true, 1, $false, "27+1=", 2) | Out-Null
$d.Content.Find.Execute("54-3=", $true, $false, $false, $false, $false, $true, 1, $false, "66+23=", 2) | Out-Null
$d.Content.Find.Execute("75-55=", $true, $false, $false, $false, $false, $true, 1, $false, "58-52=", 2) | Out-Null
$d.Content.Find.Execute("99-56=", $true, $false, $false, $false, $false, $true, 1, $false, "3+14=", 2) | Out-Null
$d.Content.Find.Execute("50-22=", $true, $false, $false, $false, $false, $true, 1, $false, "58-55=", 2) | Out-Null
$d.Content.Find.Execute("20+71=", $true, $false, $false, $false, $false, $true, 1, $false, "89-10=", 2) | Out-Null
$d.Content.Find.Execute("92-48=", $true, $false, $false, $false, $false, $true, 1, $false, "13+82=", 2) | Out-Null
$d.Content.Find.Execute("30+56=", $true, $false, $false, $false, $false, $true, 1, $false, "43+7=", 2) | Out-Null
$d.Content.Find.Execute("77-42=", $true, $false, $false, $false, $false, $true, 1, $false, "14+5=", 2) | Out-Null
$d.Content.Find.Execute("10+76=", $true, $false, $false, $false, $false, $true, 1, $false, "67+22=", 2) | Out-Null
$d.Content.Find.Execute("91-90=", $true, $false, $false, $false, $false, $true, 1, $false, "76-14=", 2) | Out-Null
$d.Content.Find.Execute("53+24=", $true, $false, $false, $false, $false, $true, 1, $false, "76-26=", 2) | Out-Null
$d.Content.Find.Execute("44+7=", $true, $false, $false, $false, $false, $true, 1, $false, "98-13=", 2) | Out-Null
$d.Content.Find.Execute("16+72=", $true, $false, $false, $false, $false, $true, 1, $false, "86+8=", 2) | Out-Null
$d.Content.Find.Execute("64+5=", $true, $false, $false, $false, $false, $true, 1, $false, "76-73=", 2) | Out-Null
$d.Content.Find.Execute("76-27=", $true, $false, $false, $false, $false, $true, 1, $false, "63+36=", 2) | Out-Null
$d.Content.Find.Execute("66-30=", $true, $false, $false, $false, $false, $true, 1, $false, "63+5=", 2) | Out-Null
$d.Content.Find.Execute("90+5=", $true, $false, $false, $false, $false, $true, 1, $false, "1+19=", 2) | Out-Null
$d.Content.Find.Execute("11+27=", $true, $false, $false, $false, $false, $true, 1, $false, "10+12=", 2) | Out-Null
$d.Content.Find.Execute("56+34=", $true, $false, $false, $false, $false, $true, 1, $false, "0+37=", 2) | Out-Null
$d.Content.Find.Execute("17+75=", $true, $false, $false, $false, $false, $true, 1, $false, "83-40=", 2) | Out-Null
$d.Content.Find.Execute("4+16=", $true, $false, $false, $false, $false, $true, 1, $false, "92-62=", 2) | Out-Null
$d.Content.Find.Execute("8+78=", $true, $false, $false, $false, $false, $true, 1, $false, "36+10=", 2) | Out-Null
$d.Content.Find.Execute("22-1=", $true, $false, $false, $false, $false, $true, 1, $false, "43+48=", 2) | Out-Null
$d.Content.Find.Execute("82-18=", $true, $false, $false, $false, $false, $true, 1, $false, "59+13=", 2) | Out-Null
$d.Content.Find.Execute("97-72=", $true, $false, $false, $false, $false, $true, 1, $false, "10+41=", 2) | Out-Null
$d.Content.Find.Execute("95-13=", $true, $false, $false, $false, $false, $true, 1, $false, "40+10=", 2) | Out-Null
$d.Content.Find.Execute("69+16=", $true, $false, $false, $false, $false, $true, 1, $false, "24-23=", 2) | Out-Null
$d.Content.Find.Execute("47+32=", $true, $false, $false, $false, $false, $true, 1, $false, "32+37=", 2) | Out-Null
$d.Content.Find.Execute("15+16=", $true, $false, $false, $false, $false, $true, 1, $false, "83+15=", 2) | Out-Null
$d.Content.Find.Execute("58-1=", $true, $false, $false, $false, $false, $true, 1, $false, "5+19=", 2) | Out-Null
$d.Content.Find.Execute("25-3=", $true, $false, $false, $false, $false, $true, 1, $false, "97-64=", 2) | Out-Null
$d.Content.Find.Execute("84-67=", $true, $false, $false, $false, $false, $true, 1, $false, "30+55=", 2) | Out-Null
$d.Content.Find.Execute("7+4=", $true, $false, $false, $false, $false, $true, 1, $false, "68+21=", 2) | Out-Null
$d.Content.Find.Execute("74-53=", $true, $false, $false, $false, $false, $true, 1, $false, "10+43=", 2) | Out-Null
$d.Content.Find.Execute("3+11=", $true, $false, $false, $false, $false, $true, 1, $false, "47-10=", 2) | Out-Null
$d.Content.Find.Execute("26-16=", $true, $false, $false, $false, $false, $true, 1, $false, "90-87=", 2) | Out-Null
$d.Content.Find.Execute("79-64=", $true, $false, $false, $false, $false, $true, 1, $false, "54+12=", 2) | Out-Null
$d.Content.Find.Execute("50-39=", $true, $false, $false, $false, $false, $true, 1, $false, "19+51=", 2) | Out-Null
$d.Content.Find.Execute("67-4=", $true, $false, $false, $false, $false, $true, 1, $false, "65+0=", 2) | Out-Null
$d.Content.Find.Execute("50-5=", $true, $false, $false, $false, $false, $true, 1, $false, "62+1=", 2) | Out-Null
$d.Content.Find.Execute("73+7=", $true, $false, $false, $false, $false, $true, 1, $false, "46-32=", 2) | Out-Null
$d.Content.Find.Execute("68-41=", $true, $false, $false, $false, $false, $true, 1, $false, "57-10=", 2) | Out-Null
$d.Content.Find.Execute("93-43=", $true, $false, $false, $false, $false, $true, 1, $false, "8+61=", 2) | Out-Null
$d.Content.Find.Execute("87-22=", $true, $false, $false, $false, $false, $true, 1, $false, "66+30=", 2) | Out-Null
$d.Content.Find.Execute("88-33=", $true, $false, $false, $false, $false, $true, 1, $false, "77-63=", 2) | Out-Null
$d.Content.Find.Execute("75-17=", $true, $false, $false, $false, $false, $true, 1, $false, "27+64=", 2) | Out-Null
$d.Content.Find.Execute("9+46=", $true, $false, $false, $false, $false, $true, 1, $false, "55-13=", 2) | Out-Null
$d.Content.Find.Execute("40-26=", $true, $false, $false, $false, $false, $true, 1, $false, "84-81=", 2) | Out-Null
$d.Content.Find.Execute("12+48=", $true, $false, $false, $false, $false, $true, 1, $false, "53-37=", 2) | Out-Null
$d.Content.Find.Execute("16+36=", $true, $false, $false, $false, $false, $true, 1, $false, "85+14=", 2) | Out-Null
$d.Content.Find.Execute("90-4=", $true, $false, $false, $false, $false, $true, 1, $false, "31+20=", 2) | Out-Null
$d.Content.Find.Execute("91-33=", $true, $false, $false, $false, $false, $true, 1, $false, "90-3=", 2) | Out-Null
$d.Content.Find.Execute("3+55=", $true, $false, $false, $false, $false, $true, 1, $false, "91-58=", 2) | Out-Null
$d.Content.Find.Execute("93-13=", $true, $false, $false, $false, $false, $true, 1, $false, "72-53=", 2) | Out-Null
$d.Content.Find.Execute("80-30=", $true, $false, $false, $false, $false, $true, 1, $false, "53-13=", 2) | Out-Null
$d.Content.Find.Execute("65-1=", $true, $false, $false, $false, $false, $true, 1, $false, "87-76=", 2) | Out-Null
$d.Content.Find.Execute("25+4=", $true, $false, $false, $false, $false, $true, 1, $false, "83-26=", 2) | Out-Null
$d.Content.Find.Execute("80+19=", $true, $false, $false, $false, $false, $true, 1, $false, "35+19=", 2) | Out-Null
$d.Content.Find.Execute("14+21=", $true, $false, $false, $false, $false, $true, 1, $false, "89-84=", 2) | Out-Null
$d.Content.Find.Execute("99-40=", $true, $false, $false, $false, $false, $true, 1, $false, "19+75=", 2) | Out-Null
$d.Content.Find.Execute("50-16=", $true, $false, $false, $false, $false, $true, 1, $false, "8+30=", 2) | Out-Null
$d.Content.Find.Execute("67-17=", $true, $false, $false, $false, $false, $true, 1, $false, "84+5=", 2) | Out-Null
$d.Content.Find.Execute("98-55=", $true, $false, $false, $false, $false, $true, 1, $false, "39-27=", 2) | Out-Null
$d.Content.Find.Execute("31+8=", $true, $false, $false, $false, $false, $true, 1, $false, "92-22=", 2) | Out-Null
